$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the Eurostat source hyperlink in F3, then split the old D3 text ---
# --- into D3 (title) + E3 (dataset code) ---
$ws.Hyperlinks.Add($ws.Range("F3"), "https://ec.europa.eu/eurostat/databrowser/view/demo_r_d2jan/default/table") | Out-Null
$ws.Range("D3").Value = "Population on 1 January by age, sex and NUTS 2 region "
$ws.Range("E3").Value = "[DEMO_R_D2JAN]"

# --- "acessed" / "last update" block (F4:G5) ---
$ws.Range("F5").Value = "acessed"
$ws.Range("G5").Value = "29.09.2023"

$ws.Range("F4").Value = "last update"
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"
$ws.Range("G4").Value = 45197.958333333336

# --- Eurostat metadata block (F6:G8) ---
$ws.Range("F6").Value = "Unit of measure [UNIT]"
$ws.Range("G6").Value = "Number [NR]"

$ws.Range("F7").Value = "Sex [SEX]"
$ws.Range("G7").Value = "Total [T]"

$ws.Range("F8").Value = "Age class [AGE]"
$ws.Range("G8").Value = "Total [TOTAL]"

# --- second source (World Data Bank / Bosnia) block, row 9 ---
$ws.Range("B9").Value = "[2]"
$ws.Range("C9").Value = "Bosnia"
$ws.Range("D9").Value = "World Data Bank, see Popular Indicators: Population total"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://databank.worldbank.org/source/population-estimates-and-projections") | Out-Null
# F10:F11 inherit the hyperlink-cell formatting (same as the source file) without a value
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F10:F11").PasteSpecial(-4122) | Out-Null

# --- exported-data stamp (D10:E11), Arial 10 font, dd.mm.yy date format ---
$ws.Range("D10").Value = "Exportierte Daten"
$ws.Range("E10").Value = "16.08.2023"
$ws.Range("D10:E10").Font.Name = "Arial"
$ws.Range("D10:E10").Font.Size = 10
$ws.Range("E10").NumberFormat = "dd.mm.yy"
$ws.Range("E10").VerticalAlignment = -4108

# D11/E11 inherit the same formatting, no values
$ws.Range("D10:E10").Copy() | Out-Null
$ws.Range("D11:E11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- column widths ---
$ws.Columns.Item(1).ColumnWidth = 11.25
$ws.Columns.Item(4).ColumnWidth = 16.751
$ws.Columns.Item(5).ColumnWidth = 15.917
$ws.Columns.Item(6).ColumnWidth = 11.1
$ws.Columns.Item(7).ColumnWidth = 16.584

# --- selection matches the saved view in the target file ---
$ws.Range("C16").Select() | Out-Null
